# The target paragraph is a single run made up of <w:t> segments separated
# by <w:br/> manual line breaks. [char]11 (vertical tab) is how Word's COM
# object model represents a <w:br/> in Range.Text, so we rebuild the whole
# block as one big string joined on that character and swap it in with a
# single Find/Replace across the whole paragraph (keeps it one run, with no
# stray xml:space="preserve" the way a direct Range.Text assignment would).
$d = $word.ActiveDocument
$nl = [char]11

$findText = @(
    'Obsah online kurzu "EU Data Act a jeho dopady na práci s daty" v Moodle:',
    '',
    '1. Úvodní sekce:',
    '- Přivítání studentů',
    '- Název kurzu: EU Data Act a jeho dopady na práci s daty',
    '- Garant: doc. JUDr. Jan Novák, Ph.D.',
    '- Kontaktní údaje: [email protected]',
    '- Cíle kurzu: Porozumět klíčovým ustanovením EU Data Actu, analyzovat jeho dopady a získat dovednosti v implementaci požadavků',
    '- Požadavky k absolvování: Základy práva v ICT nebo úvod do datové analýzy',
    '- Doporučená literatura',
    '- Harmonogram výuky',
    '',
    '2. Struktura kurzu:',
    '- Tématické týdny/moduly:',
    '  a) Úvod do EU Data Actu',
    '  b) Povinnosti výrobců a poskytovatelů dat',
    '  c) Práva uživatelů a třetích stran',
    '  d) Data Act vs. GDPR a další právní rámce',
    '  e) Interoperabilita a technické standardy',
    '  f) Data Act a cloudové služby',
    '  g) Povinnosti v oblasti veřejné správy',
    '  h) Implementace Data Actu v praxi',
    '  i) Případové studie',
    '',
    '3. Studijní materiály:',
    '- PDF/textové soubory s teoretickými informacemi k jednotlivým témům',
    '- Video přednášky od garantujícího doc. JUDr. Jana Nováka, Ph.D.',
    '- Prezentace ke stažení ve formátu PowerPoint nebo PDF',
    '- Odkazy na relevantní zdroje k problematice EU Data Actu',
    '',
    '4. Aktivizace studenta:',
    '- Diskuzní fóra na konci každého tématického týdne',
    '- Samotestovací kvízy k ověření znalostí',
    '- Průběžné úkoly k procvičení praktických dovedností',
    '- Projekty nebo eseje k získání praktických zkušeností',
    '',
    '5. Hodnocení a zpětná vazba:',
    '- Automaticky hodnocené testy k ověření znalostí',
    '- Ručně hodnocené úkoly k zhodnocení praktických dovedností',
    '- Závěrečný test nebo projekt k celkovému zhodnocení kurzu',
    '- Možnost studentova sebehodnocení',
    '- Možnost komentářů a zpětné vazby od garantujícího',
    '',
    '6. Zpřístupnění a dostupnost:',
    '- Postupné zpřístupnění obsahu kurzu v jednotlivých tématických týdnech',
    '- Dostupnost celého kurzu pro předstihové studium',
    '- Soubory ve formátech vhodných i pro mobilní zařízení',
    '- Popisky pro multimediální obsah pro lepší přístupnost',
    '',
    '7. Jazyk a styl:',
    '- Formální, ale srozumitelný jazyk v souladu s cílovou skupinou magisterského studia',
    '- Strukturované a přehledné prezentace informací',
    '- Zvýraznění klíčových bodů pro snadnější orientaci',
    '',
    '8. Závěrečná sekce kurzu:',
    '- Shrnutí klíčových bodů a poznatků z kurzu',
    '- Závěrečný test nebo projekt k ukončení kurzu',
    '- Možnost studentů vyjádřit zpětnou vazbu a poděkování',
    '- Informace o dalších kurzech nebo tématech navazujících na problematiku EU Data Actu',
    '',
    '9. Technické doporučení:',
    '- Standardizované názvy souborů a sekcí pro usnadnění orientace',
    '- Omezení velikosti souborů pro rychlejší načítání',
    '- Pravidelné zálohování obsahu kurzu pro zajištění bezpečnosti dat',
    '- Testování kurzu z pohledu studenta pro odhalení případných technických nedostatků',
    '',
    'Tímto způsobem by měl být online kurz EU Data Actu strukturován a připraven pro efektivní vzdělávání studentů magisterského studia.'
) -join $nl

$replText = @(
    'The course you are tasked with creating is called ¦Data Act Course.¦ Based on the provided methodology, the course should be structured with the following sections in Moodle:',
    '',
    '1. **Announcements**: This section is for course communication with students during the semester. It should include organizational updates and important announcements.',
    '',
    '2. **Introduction**: Provide basic course information such as course focus, instructors, learning outcomes, grading method, a link to the official syllabus, general rules, and course structure. You may also include a short introductory survey or quiz.',
    '',
    '3. **Modules**: Propose 4–8 logically grouped modules based on the syllabus. Each module should have a numbered title, introductory paragraph, structured content explanation, key point summary, and a self-check quiz.',
    '',
    '4. **Conclusion**: Include information about how the course ends, the final exam or project, and space for student feedback.',
    '',
    'Ensure consistency in tone, use Markdown headers, bullet points for lists, and follow the provided Quiz Format Requirements for self-check quizzes and the final comprehensive quiz.',
    '',
    'If you have any specific questions or need further assistance, feel free to ask!'
) -join $nl

$rng = $d.Content
$found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replText, 2)
if (-not $found) {
    throw "Could not locate the course-outline paragraph to replace."
}

# Find/Replace auto-"smart-quotes" any literal " in the replacement text, but
# the target text needs straight quotes, so we swapped them for a placeholder
# character above and fix each one back up individually afterwards (a small
# Range.Text assignment on just that character keeps straight quotes).
$placeholder = [char]166
$fixupRange = $d.Content
while ($fixupRange.Find.Execute($placeholder, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)) {
    $fixupRange.Text = '"'
    $fixupRange = $d.Range($fixupRange.End, $d.Content.End)
}
